# Apply targeted text corrections/updates to worksheet "Tab_5a_Indikatoren"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_5a_Indikatoren")

# --- Row 21: hyphenate "17 und 18-Jährigen" -> "17- und 18-Jährigen" ---
$ws.Range("D21").Value = "Anteil der 17- und 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Range("F21").Value = "Anteil der 17- und 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Range("L21").Value = "Anteil der 17- und 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Range("M21").Value = "XXXAnteil der 17- und 18-Jährigen mit (angestrebter) Studienberechtigung"

# --- Row 70: climate finance target wording update ---
$ws.Range("H70").Value = "Erhöhung auf jährlich mindestens 6 Milliarden Euro bis spätestens 2025"
$ws.Range("I70").Value = "XXXErhöhung auf jährlich mindestens 6 Milliarden Euro bis spätestens 2025"

# --- Row 73: fish stock management target wording + year update ---
$ws.Range("H73").Value = "Anteil nachhaltig bewirtschafteter Fischbestände in Nord- und Ostsee an allen MSY-untersuchten Beständen soll bis 2030 100 Prozent betragen"
$ws.Range("I73").Value = "XXXAnteil nachhaltig bewirtschafteter Fischbestände in Nord- und Ostsee an allen MSY-untersuchten Beständen soll bis 2030 100 Prozent betragen"
$ws.Range("J73").Value = "nachhaltige Bewirtschaftung nach dem MSY-Ansatz bis 2030"
$ws.Range("K73").Value = "sustainable management in accordance with the MSY approachby 2030"

# --- Row 76: replace REDD+ indicator with Bodenversiegelungsgrad (soil sealing) ---
$ws.Range("A76").Value = "Z15_B03_P01_Ib01_I01"
$ws.Range("B76").Value = "Z15_B03_P01_Ib01"
$ws.Range("C76").Value = "15.3"
$ws.Range("D76").Value = "Bodenversiegelungsgrad"
$ws.Range("E76").Value = "XXXBodenversiegelungsgrad"
$ws.Range("F76").Value = "Bodenversiegelungsgrad"
$ws.Range("G76").Value = "XXXBodenversiegelungsgrad"
$ws.Range("H76").Value = "Sinkende Zunahme der Bodenversiegelung"
$ws.Range("I76").Value = "XXXSinkende Zunahme der Bodenversiegelung"
$ws.Range("J76").Value = "Sinkende Zunahme"
$ws.Range("K76").Value = "XXXSinkende Zunahme"
$ws.Range("L76").Value = "Bodenversiegelungsgrad"
$ws.Range("M76").Value = "XXXBodenversiegelungsgrad"

# --- Row 77: replace soil-protection indicator with bilateral land contributions ---
$ws.Range("A77").Value = "Z15_B04_P01_Ib01_I01"
$ws.Range("B77").Value = "Z15_B04_P01_Ib01"
$ws.Range("C77").Value = "15.4"
$ws.Range("D77").Value = "Bilaterale Beiträge der deutschen internationalen Kooperation zum Schutz, nachhaltiger Nutzung und Wiederherstellung von Land, gemessen in Hektar"
$ws.Range("E77").Value = "XXXBilaterale Beiträge der deutschen internationalen Kooperation zum Schutz, nachhaltiger Nutzung und Wiederherstellung von Land, gemessen in Hektar"
$ws.Range("F77").Value = "Bilaterale Beiträge der deutschen internationalen Kooperation zum Schutz, nachhaltiger Nutzung und Wiederherstellung von Land, gemessen in Hektar"
$ws.Range("G77").Value = "XXXBilaterale Beiträge der deutschen internationalen Kooperation zum Schutz, nachhaltiger Nutzung und Wiederherstellung von Land, gemessen in Hektar"
$ws.Range("H77").Value = "Steigerung der nachhaltigen Landnutzung (durch Schutz, nachhaltige Bewirtschaftung, Wiederherstellung) bis 2030"
$ws.Range("I77").Value = "XXXSteigerung der nachhaltigen Landnutzung (durch Schutz, nachhaltige Bewirtschaftung, Wiederherstellung) bis 2030"
$ws.Range("J77").Value = "Steigerung der nachhaltigen Landnutzung"
$ws.Range("K77").Value = "XXXSteigerung der nachhaltigen Landnutzung"
$ws.Range("L77").Value = "Bilaterale Beiträge der deutschen internationalen Kooperation zum Schutz, nachhaltiger Nutzung und Wiederherstellung von Land, gemessen in Hektar"
$ws.Range("M77").Value = "XXXBilaterale Beiträge der deutschen internationalen Kooperation zum Schutz, nachhaltiger Nutzung und Wiederherstellung von Land, gemessen in Hektar"
